# The commit deletes the "Course Content" slide (slide id 296), which was the
# 5th slide in the deck, a leftover placeholder slide containing dummy text
# ("Dfhkalsjdf", "Asdf", "Fasdhfsad", "sadf"). Removing it shifts every slide
# after it up by one position; PowerPoint/the host take care of renumbering
# the remaining slide IDs, relationship IDs, cached slide-number fields, and
# section membership automatically when the slide is deleted.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$s.Delete()
